$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - Esperanza vida (fill in remaining indicator data)
$ws.Range("C7").Value = 1993
$ws.Range("D7").Value = 2018
$ws.Range("E7").Value = 2019
$ws.Range("F7").Value = 2031
$ws.Range("G7").Value = 1993
$ws.Range("H7").Value = 2031
$ws.Range("I7").Value = "SI"
$ws.Range("J7").Value = "SI"
$ws.Range("K7").Value = "NO"
$ws.Range("L7").Value = "NO"

# Row 8 - Indicador 07 -> Crímenes
$ws.Range("B8").Value = "Crímenes"
$ws.Range("C8").Value = 2010
$ws.Range("D8").Value = 2022
$ws.Range("E8").Value = 2023
$ws.Range("F8").Value = 2031
$ws.Range("G8").Value = 2010
$ws.Range("H8").Value = 2031
$ws.Range("I8").Value = "NO / SUMA"
$ws.Range("J8").Value = "NO"
$ws.Range("K8").Value = "NO"
$ws.Range("L8").Value = "NO"

# Row 9 - Indicador 08 -> Servicios
$ws.Range("B9").Value = "Servicios"
$ws.Range("C9").Value = 2009
$ws.Range("D9").Value = 2024
$ws.Range("E9").Value = 2025
$ws.Range("F9").Value = 2031
$ws.Range("G9").Value = 2009
$ws.Range("H9").Value = 2031
$ws.Range("I9").Value = "SI"
$ws.Range("J9").Value = "SI"
$ws.Range("K9").Value = "SI"
$ws.Range("L9").Value = "NO"

# Update selection to match final cursor position
$ws.Range("L10").Select()
